$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1308.8572  # H6: 1021.8889 -> 1308.8572
$ws.Cells.Item(6, 9).Value = 693.6667  # I6: 524.625 -> 693.6667
$ws.Cells.Item(6, 11).Value = 2081.0001  # K6: 1573.875 -> 2081.0001
$ws.Cells.Item(6, 13).Value = -1969.0001  # M6: -1461.875 -> -1969.0001
$ws.Cells.Item(8, 8).Value = 59.625  # H8: 67.42856999999999 -> 59.625
$ws.Cells.Item(8, 9).Value = 59.625  # I8: 67.42856999999999 -> 59.625
$ws.Cells.Item(8, 11).Value = 178.875  # K8: 202.28571 -> 178.875
$ws.Cells.Item(8, 13).Value = -39.875  # M8: -63.28570999999999 -> -39.875
$ws.Cells.Item(33, 8).Value = 338.14285  # H33: 340.85715 -> 338.14285
$ws.Cells.Item(33, 9).Value = 244.58333  # I33: 247.75 -> 244.58333
$ws.Cells.Item(33, 11).Value = 244.58333  # K33: 247.75 -> 244.58333
$ws.Cells.Item(33, 13).Value = -15.58332999999999  # M33: -18.75 -> -15.58332999999999
$ws.Cells.Item(51, 8).Value = 6634.6  # H51: 6714.6 -> 6634.6
$ws.Cells.Item(51, 9).Value = 8196.666999999999  # I51: 6747.5 -> 8196.666999999999
$ws.Cells.Item(51, 10).Value = 4291.5  # J51: 6583 -> 4291.5
$ws.Cells.Item(51, 11).Value = 8196.666999999999  # K51: 6747.5 -> 8196.666999999999
$ws.Cells.Item(51, 12).Value = 4291.5  # L51: 6583 -> 4291.5
$ws.Cells.Item(51, 13).Value = -7712.666999999999  # M51: -6263.5 -> -7712.666999999999
$ws.Cells.Item(51, 14).Value = -5259.5  # N51: -7551 -> -5259.5
$ws.Cells.Item(52, 8).Value = 709  # H52: 400 -> 709
$ws.Cells.Item(52, 9).Value = 709  # I52: 0 -> 709
$ws.Cells.Item(52, 10).Value = 0  # J52: 400 -> 0
$ws.Cells.Item(52, 11).Value = 2127  # K52: 0 -> 2127
$ws.Cells.Item(52, 12).Value = 0  # L52: 1200 -> 0
$ws.Cells.Item(52, 13).Value = -1967  # M52: <<ABSENT>> -> -1967
$ws.Cells.Item(52, 14).ClearContents()  # N52: -1520 -> (removed)
$ws.Cells.Item(64, 8).Value = 7575.4614  # H64: 8036.846 -> 7575.4614
$ws.Cells.Item(64, 9).Value = 3700.2  # I64: 4624.75 -> 3700.2
$ws.Cells.Item(64, 10).Value = 9997.5  # J64: 9553.333000000001 -> 9997.5
$ws.Cells.Item(64, 11).Value = 3700.2  # K64: 4624.75 -> 3700.2
$ws.Cells.Item(64, 12).Value = 9997.5  # L64: 9553.333000000001 -> 9997.5
$ws.Cells.Item(64, 13).Value = -3452.2  # M64: -4376.75 -> -3452.2
$ws.Cells.Item(64, 14).Value = -10493.5  # N64: -10049.333 -> -10493.5
$ws.Cells.Item(67, 8).Value = 7575.4614  # H67: 8036.846 -> 7575.4614
$ws.Cells.Item(67, 9).Value = 3700.2  # I67: 4624.75 -> 3700.2
$ws.Cells.Item(67, 10).Value = 9997.5  # J67: 9553.333000000001 -> 9997.5
$ws.Cells.Item(67, 11).Value = 3700.2  # K67: 4624.75 -> 3700.2
$ws.Cells.Item(67, 12).Value = 9997.5  # L67: 9553.333000000001 -> 9997.5
$ws.Cells.Item(67, 13).Value = -2842.2  # M67: -3766.75 -> -2842.2
$ws.Cells.Item(67, 14).Value = -11713.5  # N67: -11269.333 -> -11713.5
$ws.Cells.Item(70, 8).Value = 39581.6  # H70: 18832.273 -> 39581.6
$ws.Cells.Item(70, 9).Value = 3350  # I70: 1200 -> 3350
$ws.Cells.Item(70, 10).Value = 48639.5  # J70: 20595.5 -> 48639.5
$ws.Cells.Item(70, 11).Value = 10050  # K70: 3600 -> 10050
$ws.Cells.Item(70, 12).Value = 145918.5  # L70: 61786.5 -> 145918.5
$ws.Cells.Item(70, 13).Value = -9780  # M70: -3330 -> -9780
$ws.Cells.Item(70, 14).Value = -146458.5  # N70: -62326.5 -> -146458.5
$ws.Cells.Item(73, 8).Value = 39581.6  # H73: 18832.273 -> 39581.6
$ws.Cells.Item(73, 9).Value = 3350  # I73: 1200 -> 3350
$ws.Cells.Item(73, 10).Value = 48639.5  # J73: 20595.5 -> 48639.5
$ws.Cells.Item(73, 11).Value = 10050  # K73: 3600 -> 10050
$ws.Cells.Item(73, 12).Value = 145918.5  # L73: 61786.5 -> 145918.5
$ws.Cells.Item(73, 13).Value = -9114  # M73: -2664 -> -9114
$ws.Cells.Item(73, 14).Value = -147790.5  # N73: -63658.5 -> -147790.5
$ws.Cells.Item(74, 8).Value = 8302.958000000001  # H74: 8305.083000000001 -> 8302.958000000001
$ws.Cells.Item(74, 9).Value = 5549.6665  # I74: 5566.3335 -> 5549.6665
$ws.Cells.Item(74, 10).Value = 8696.286  # J74: 8696.333000000001 -> 8696.286
$ws.Cells.Item(74, 11).Value = 5549.6665  # K74: 5566.3335 -> 5549.6665
$ws.Cells.Item(74, 12).Value = 8696.286  # L74: 8696.333000000001 -> 8696.286
$ws.Cells.Item(74, 13).Value = -4613.6665  # M74: -4630.3335 -> -4613.6665
$ws.Cells.Item(74, 14).Value = -10568.286  # N74: -10568.333 -> -10568.286
$ws.Cells.Item(77, 8).Value = 8302.958000000001  # H77: 8305.083000000001 -> 8302.958000000001
$ws.Cells.Item(77, 9).Value = 5549.6665  # I77: 5566.3335 -> 5549.6665
$ws.Cells.Item(77, 10).Value = 8696.286  # J77: 8696.333000000001 -> 8696.286
$ws.Cells.Item(77, 11).Value = 27748.3325  # K77: 27831.6675 -> 27748.3325
$ws.Cells.Item(77, 12).Value = 43481.43  # L77: 43481.665 -> 43481.43
$ws.Cells.Item(77, 13).Value = -23068.3325  # M77: -23151.6675 -> -23068.3325
$ws.Cells.Item(77, 14).Value = -52841.43  # N77: -52841.665 -> -52841.43
$ws.Cells.Item(112, 8).Value = 953.29034  # H112: 958.4138 -> 953.29034
$ws.Cells.Item(112, 10).Value = 929.37933  # J112: 933.1111 -> 929.37933
$ws.Cells.Item(112, 12).Value = 2788.13799  # L112: 2799.3333 -> 2788.13799
$ws.Cells.Item(112, 14).Value = -5004.13799  # N112: -5015.3333 -> -5004.13799
$ws.Cells.Item(116, 8).Value = 111331.836  # H116: 76533.44500000001 -> 111331.836
$ws.Cells.Item(116, 10).Value = 35330  # J116: 21133.334 -> 35330
$ws.Cells.Item(116, 12).Value = 35330  # L116: 21133.334 -> 35330
$ws.Cells.Item(116, 14).Value = -42214  # N116: -28017.334 -> -42214
$ws.Cells.Item(137, 8).Value = 31252484  # H137: 26318072 -> 31252484
$ws.Cells.Item(137, 10).Value = 2631.2856  # J137: 2201.3 -> 2631.2856
$ws.Cells.Item(137, 12).Value = 7893.8568  # L137: 6603.900000000001 -> 7893.8568
$ws.Cells.Item(137, 14).Value = -12993.8568  # N137: -11703.9 -> -12993.8568
$ws.Cells.Item(138, 8).Value = 747.52  # H138: 749.4286 -> 747.52
$ws.Cells.Item(138, 9).Value = 747.52  # I138: 749.4286 -> 747.52
$ws.Cells.Item(138, 11).Value = 2242.56  # K138: 2248.2858 -> 2242.56
$ws.Cells.Item(138, 13).Value = 2897.44  # M138: 2891.7142 -> 2897.44

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1994.9242  # H32: 2112.258 -> 1994.9242
$ws.Cells.Item(32, 9).Value = 1142.0702  # I32: 1214.9623 -> 1142.0702
$ws.Cells.Item(32, 11).Value = 1142.0702  # K32: 1214.9623 -> 1142.0702
$ws.Cells.Item(32, 13).Value = -855.0702000000001  # M32: -927.9622999999999 -> -855.0702000000001
$ws.Cells.Item(45, 8).Value = 2714.4443  # H45: 2448.353 -> 2714.4443
$ws.Cells.Item(45, 9).Value = 2554  # I45: 2187.5715 -> 2554
$ws.Cells.Item(45, 10).Value = 3998  # J45: 3665.3333 -> 3998
$ws.Cells.Item(45, 11).Value = 2554  # K45: 2187.5715 -> 2554
$ws.Cells.Item(45, 12).Value = 3998  # L45: 3665.3333 -> 3998
$ws.Cells.Item(45, 13).Value = -2177  # M45: -1810.5715 -> -2177
$ws.Cells.Item(45, 14).Value = -4752  # N45: -4419.3333 -> -4752
$ws.Cells.Item(61, 8).Value = 1992  # H61: 1960.081 -> 1992
$ws.Cells.Item(61, 9).Value = 1829.5862  # I61: 1763.5161 -> 1829.5862
$ws.Cells.Item(61, 10).Value = 2664.8572  # J61: 2975.6667 -> 2664.8572
$ws.Cells.Item(61, 11).Value = 1829.5862  # K61: 1763.5161 -> 1829.5862
$ws.Cells.Item(61, 12).Value = 2664.8572  # L61: 2975.6667 -> 2664.8572
$ws.Cells.Item(61, 13).Value = -1617.5862  # M61: -1551.5161 -> -1617.5862
$ws.Cells.Item(61, 14).Value = -3088.8572  # N61: -3399.6667 -> -3088.8572
$ws.Cells.Item(122, 8).Value = 1552.88  # H122: 1543.1538 -> 1552.88
$ws.Cells.Item(122, 9).Value = 1491.7273  # I122: 1483.3914 -> 1491.7273
$ws.Cells.Item(122, 11).Value = 4475.1819  # K122: 4450.174199999999 -> 4475.1819
$ws.Cells.Item(122, 13).Value = -2025.1819  # M122: -2000.174199999999 -> -2025.1819
$ws.Cells.Item(123, 8).Value = 0  # H123: 68000 -> 0
$ws.Cells.Item(123, 10).Value = 0  # J123: 68000 -> 0
$ws.Cells.Item(123, 12).Value = 0  # L123: 68000 -> 0
$ws.Cells.Item(123, 14).ClearContents()  # N123: -77800 -> (removed)
$ws.Cells.Item(132, 8).Value = 4606.5386  # H132: 5087.222 -> 4606.5386
$ws.Cells.Item(132, 9).Value = 4652.364  # I132: 5296.5713 -> 4652.364
$ws.Cells.Item(132, 11).Value = 13957.092  # K132: 15889.7139 -> 13957.092
$ws.Cells.Item(132, 13).Value = -11427.092  # M132: -13359.7139 -> -11427.092
$ws.Cells.Item(133, 8).Value = 99000  # H133: 98984.5 -> 99000
$ws.Cells.Item(133, 10).Value = 99000  # J133: 98984.5 -> 99000
$ws.Cells.Item(133, 12).Value = 99000  # L133: 98984.5 -> 99000
$ws.Cells.Item(133, 14).Value = -104060  # N133: -104044.5 -> -104060
$ws.Cells.Item(136, 8).Value = 1992  # H136: 1960.081 -> 1992
$ws.Cells.Item(136, 9).Value = 1829.5862  # I136: 1763.5161 -> 1829.5862
$ws.Cells.Item(136, 10).Value = 2664.8572  # J136: 2975.6667 -> 2664.8572
$ws.Cells.Item(136, 11).Value = 5488.7586  # K136: 5290.5483 -> 5488.7586
$ws.Cells.Item(136, 12).Value = 7994.571599999999  # L136: 8927.000100000001 -> 7994.571599999999
$ws.Cells.Item(136, 13).Value = -2938.7586  # M136: -2740.5483 -> -2938.7586
$ws.Cells.Item(136, 14).Value = -13094.5716  # N136: -14027.0001 -> -13094.5716

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1881.2  # H86: 2051.5 -> 1881.2
$ws.Cells.Item(86, 9).Value = 1601.5  # I86: 1735.3334 -> 1601.5
$ws.Cells.Item(86, 11).Value = 1601.5  # K86: 1735.3334 -> 1601.5
$ws.Cells.Item(86, 13).Value = -478.5  # M86: -612.3334 -> -478.5
$ws.Cells.Item(89, 8).Value = 1881.2  # H89: 2051.5 -> 1881.2
$ws.Cells.Item(89, 9).Value = 1601.5  # I89: 1735.3334 -> 1601.5
$ws.Cells.Item(89, 11).Value = 8007.5  # K89: 8676.666999999999 -> 8007.5
$ws.Cells.Item(89, 13).Value = -2391.5  # M89: -3060.666999999999 -> -2391.5
$ws.Cells.Item(94, 8).Value = 975.2  # H94: 893.7857 -> 975.2
$ws.Cells.Item(94, 9).Value = 957.1579  # I94: 875.381 -> 957.1579
$ws.Cells.Item(94, 10).Value = 1032.3334  # J94: 949 -> 1032.3334
$ws.Cells.Item(94, 11).Value = 957.1579  # K94: 875.381 -> 957.1579
$ws.Cells.Item(94, 12).Value = 1032.3334  # L94: 949 -> 1032.3334
$ws.Cells.Item(94, 13).Value = -506.1579  # M94: -424.381 -> -506.1579
$ws.Cells.Item(94, 14).Value = -1934.3334  # N94: -1851 -> -1934.3334
$ws.Cells.Item(107, 8).Value = 33340502  # H107: 20838278 -> 33340502
$ws.Cells.Item(107, 9).Value = 8647.25  # I107: 6050.25 -> 8647.25
$ws.Cells.Item(107, 10).Value = 166667920  # J107: 83334960 -> 166667920
$ws.Cells.Item(107, 11).Value = 8647.25  # K107: 6050.25 -> 8647.25
$ws.Cells.Item(107, 12).Value = 166667920  # L107: 83334960 -> 166667920
$ws.Cells.Item(107, 13).Value = -6727.25  # M107: -4130.25 -> -6727.25
$ws.Cells.Item(107, 14).Value = -166671760  # N107: -83338800 -> -166671760

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5625.804  # H31: 5527.1704 -> 5625.804
$ws.Cells.Item(31, 9).Value = 5842.857  # I31: 5622.273 -> 5842.857
$ws.Cells.Item(31, 11).Value = 5842.857  # K31: 5622.273 -> 5842.857
$ws.Cells.Item(31, 13).Value = -5547.857  # M31: -5327.273 -> -5547.857
$ws.Cells.Item(34, 8).Value = 5625.804  # H34: 5527.1704 -> 5625.804
$ws.Cells.Item(34, 9).Value = 5842.857  # I34: 5622.273 -> 5842.857
$ws.Cells.Item(34, 11).Value = 5842.857  # K34: 5622.273 -> 5842.857
$ws.Cells.Item(34, 13).Value = -5640.857  # M34: -5420.273 -> -5640.857
$ws.Cells.Item(58, 9).Value = 1208.6774  # I58: 1256.6774 -> 1208.6774
$ws.Cells.Item(58, 10).Value = 4450.3  # J58: 4146.727 -> 4450.3
$ws.Cells.Item(58, 11).Value = 1208.6774  # K58: 1256.6774 -> 1208.6774
$ws.Cells.Item(58, 12).Value = 4450.3  # L58: 4146.727 -> 4450.3
$ws.Cells.Item(58, 13).Value = -1005.6774  # M58: -1053.6774 -> -1005.6774
$ws.Cells.Item(58, 14).Value = -4856.3  # N58: -4552.727 -> -4856.3
$ws.Cells.Item(132, 8).Value = 166668160  # H132: 71430160 -> 166668160
$ws.Cells.Item(132, 9).Value = 166668160  # I132: 71430160 -> 166668160
$ws.Cells.Item(132, 11).Value = 500004480  # K132: 214290480 -> 500004480
$ws.Cells.Item(132, 13).Value = -500001950  # M132: -214287950 -> -500001950
$ws.Cells.Item(136, 8).Value = 1999.317  # H136: 2013.5952 -> 1999.317
$ws.Cells.Item(136, 9).Value = 1208.6774  # I136: 1256.6774 -> 1208.6774
$ws.Cells.Item(136, 10).Value = 4450.3  # J136: 4146.727 -> 4450.3
$ws.Cells.Item(136, 11).Value = 3626.0322  # K136: 3770.0322 -> 3626.0322
$ws.Cells.Item(136, 12).Value = 13350.9  # L136: 12440.181 -> 13350.9
$ws.Cells.Item(136, 13).Value = -1076.0322  # M136: -1220.0322 -> -1076.0322
$ws.Cells.Item(136, 14).Value = -18450.9  # N136: -17540.181 -> -18450.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 5995  # H80: 5447.5 -> 5995
$ws.Cells.Item(80, 10).Value = 5995  # J80: 5447.5 -> 5995
$ws.Cells.Item(80, 12).Value = 17985  # L80: 16342.5 -> 17985
$ws.Cells.Item(80, 14).Value = -19857  # N80: -18214.5 -> -19857
$ws.Cells.Item(83, 8).Value = 5995  # H83: 5447.5 -> 5995
$ws.Cells.Item(83, 10).Value = 5995  # J83: 5447.5 -> 5995
$ws.Cells.Item(83, 12).Value = 53955  # L83: 49027.5 -> 53955
$ws.Cells.Item(83, 14).Value = -63315  # N83: -58387.5 -> -63315
$ws.Cells.Item(92, 8).Value = 634.1667  # H92: 691.8 -> 634.1667
$ws.Cells.Item(92, 9).Value = 490.875  # I92: 505.2857 -> 490.875
$ws.Cells.Item(92, 10).Value = 920.75  # J92: 1127 -> 920.75
$ws.Cells.Item(92, 11).Value = 1472.625  # K92: 1515.8571 -> 1472.625
$ws.Cells.Item(92, 12).Value = 2762.25  # L92: 3381 -> 2762.25
$ws.Cells.Item(92, 13).Value = -224.625  # M92: -267.8571000000002 -> -224.625
$ws.Cells.Item(92, 14).Value = -5258.25  # N92: -5877 -> -5258.25
$ws.Cells.Item(134, 8).Value = 2193.25  # H134: 2500 -> 2193.25
$ws.Cells.Item(134, 9).Value = 2193.25  # I134: 2500 -> 2193.25
$ws.Cells.Item(134, 10).Value = 0  # J134: 2500 -> 0
$ws.Cells.Item(134, 11).Value = 6579.75  # K134: 7500 -> 6579.75
$ws.Cells.Item(134, 12).Value = 0  # L134: 7500 -> 0
$ws.Cells.Item(134, 13).Value = -1509.75  # M134: -2430 -> -1509.75
$ws.Cells.Item(134, 14).ClearContents()  # N134: -17640 -> (removed)
$ws.Cells.Item(137, 8).Value = 4167.222  # H137: 4174.593 -> 4167.222
$ws.Cells.Item(137, 10).Value = 4772.5  # J137: 4782.45 -> 4772.5
$ws.Cells.Item(137, 12).Value = 14317.5  # L137: 14347.35 -> 14317.5
$ws.Cells.Item(137, 14).Value = -24517.5  # N137: -24547.35 -> -24517.5
$ws.Cells.Item(138, 8).Value = 1793  # H138: 2525.8 -> 1793
$ws.Cells.Item(138, 9).Value = 1224  # I138: 1814.5 -> 1224
$ws.Cells.Item(138, 10).Value = 3500  # J138: 3000 -> 3500
$ws.Cells.Item(138, 11).Value = 3672  # K138: 5443.5 -> 3672
$ws.Cells.Item(138, 12).Value = 10500  # L138: 9000 -> 10500
$ws.Cells.Item(138, 13).Value = 1468  # M138: -303.5 -> 1468
$ws.Cells.Item(138, 14).Value = -20780  # N138: -19280 -> -20780

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 265.34375  # H2: 267.46875 -> 265.34375
$ws.Cells.Item(2, 9).Value = 169.57895  # I2: 177.27777 -> 169.57895
$ws.Cells.Item(2, 10).Value = 405.30768  # J2: 383.42856 -> 405.30768
$ws.Cells.Item(2, 11).Value = 169.57895  # K2: 177.27777 -> 169.57895
$ws.Cells.Item(2, 12).Value = 405.30768  # L2: 383.42856 -> 405.30768
$ws.Cells.Item(2, 13).Value = -56.57894999999999  # M2: -64.27777 -> -56.57894999999999
$ws.Cells.Item(2, 14).Value = -631.30768  # N2: -609.4285600000001 -> -631.30768

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 983  # H14: 2147.6667 -> 983
$ws.Cells.Item(14, 9).Value = 983  # I14: 2147.6667 -> 983
$ws.Cells.Item(14, 11).Value = 983  # K14: 2147.6667 -> 983
$ws.Cells.Item(14, 13).Value = -811  # M14: -1975.6667 -> -811
$ws.Cells.Item(40, 8).Value = 3051.762  # H40: 3728.5293 -> 3051.762
$ws.Cells.Item(40, 9).Value = 3063.6843  # I40: 3670.3572 -> 3063.6843
$ws.Cells.Item(40, 10).Value = 2938.5  # J40: 4000 -> 2938.5
$ws.Cells.Item(40, 11).Value = 3063.6843  # K40: 3670.3572 -> 3063.6843
$ws.Cells.Item(40, 12).Value = 2938.5  # L40: 4000 -> 2938.5
$ws.Cells.Item(40, 13).Value = -2927.6843  # M40: -3534.3572 -> -2927.6843
$ws.Cells.Item(40, 14).Value = -3210.5  # N40: -4272 -> -3210.5
$ws.Cells.Item(68, 8).Value = 3464.5  # H68: 4028 -> 3464.5
$ws.Cells.Item(68, 9).Value = 2772.182  # I68: 2833.1667 -> 2772.182
$ws.Cells.Item(68, 10).Value = 6003  # J68: 5461.8 -> 6003
$ws.Cells.Item(68, 11).Value = 2772.182  # K68: 2833.1667 -> 2772.182
$ws.Cells.Item(68, 12).Value = 6003  # L68: 5461.8 -> 6003
$ws.Cells.Item(68, 13).Value = -2023.182  # M68: -2084.1667 -> -2023.182
$ws.Cells.Item(68, 14).Value = -7501  # N68: -6959.8 -> -7501
$ws.Cells.Item(71, 8).Value = 3464.5  # H71: 4028 -> 3464.5
$ws.Cells.Item(71, 9).Value = 2772.182  # I71: 2833.1667 -> 2772.182
$ws.Cells.Item(71, 10).Value = 6003  # J71: 5461.8 -> 6003
$ws.Cells.Item(71, 11).Value = 13860.91  # K71: 14165.8335 -> 13860.91
$ws.Cells.Item(71, 12).Value = 30015  # L71: 27309 -> 30015
$ws.Cells.Item(71, 13).Value = -10116.91  # M71: -10421.8335 -> -10116.91
$ws.Cells.Item(71, 14).Value = -37503  # N71: -34797 -> -37503
$ws.Cells.Item(132, 8).Value = 2181.7083  # H132: 1914.4865 -> 2181.7083
$ws.Cells.Item(132, 9).Value = 2181.7083  # I132: 1970.138 -> 2181.7083
$ws.Cells.Item(132, 10).Value = 0  # J132: 1712.75 -> 0
$ws.Cells.Item(132, 11).Value = 6545.124899999999  # K132: 5910.414 -> 6545.124899999999
$ws.Cells.Item(132, 12).Value = 0  # L132: 5138.25 -> 0
$ws.Cells.Item(132, 13).Value = -4015.124899999999  # M132: -3380.414 -> -4015.124899999999
$ws.Cells.Item(132, 14).ClearContents()  # N132: -10198.25 -> (removed)

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1506.6428  # H122: 1367.6471 -> 1506.6428
$ws.Cells.Item(122, 9).Value = 1506.6428  # I122: 1367.6471 -> 1506.6428
$ws.Cells.Item(122, 11).Value = 4519.928400000001  # K122: 4102.9413 -> 4519.928400000001
$ws.Cells.Item(122, 13).Value = -2069.928400000001  # M122: -1652.9413 -> -2069.928400000001
